# Updates for 22 April
# Adds a new date column (AN) " 4/21/20" with that day's cumulative
# death counts for every state row, mirroring the style of the
# preceding date column (AM).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Carry over the formatting of the header cell and the regular data
# column (AL, which uses the "interior" data style rather than the
# special right-edge style reserved for the rightmost column) into the
# new column AN so that AN becomes the new right-hand edge while AM
# keeps its existing right-edge formatting.
$ws.Range("AL1:AL54").Copy()
$ws.Range("AN1:AN54").PasteSpecial(-4122)

$ws.Cells.Item(1, 40).Value = " 4/21/20"

# New daily death counts per state, in row order (rows 2-54).
$values = @(186,9,208,43,1322,486,1423,82,112,867,818,5,12,51,1468,630,83,109,171,1405,36,652,1961,2700,160,183,220,12,33,163,42,4753,65,19693,241,13,557,164,78,1614,64,171,135,8,157,528,32,40,324,682,26,242,6)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 40).Value = $values[$i]
}

# Keep the active selection pointing at the newly-added latest date
# column (Excel leaves the cursor on the just-entered cell).
$ws.Range("AN2").Select()

